$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new column E (copy the format of D1 so it matches the other headers)
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 5).PasteSpecial(-4122)
$ws.Cells.Item(1, 5).Value = "Precios"

# Fill E2:E79 with sequential numbers 1..78
$n = 1
For ($r = 2; $r -le 79; $r++) {
    $ws.Cells.Item($r, 5).Value = $n
    $n = $n + 1
}

# Update selection to G76
$ws.Range("G76").Select()
